# feat: add 2022-Q1 data
#
# 1) Create a new "2022-Q1" sheet (cloned from "2021-Q4" so it keeps the
#    exact same layout/styling), positioned between "2021-Q4" and "总计",
#    then overwrite its holding-specific figures.
# 2) Prepend a "2022-Q1" row to the "总计" (grand-total) summary sheet,
#    pushing the existing rows down and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item(4)
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(5)
$q1.Name = "2022-Q1"

# The fund code / name / rank columns are identical to 2021-Q4; only the
# size / position figures change. Force the cells to Text beforehand so
# the numeric-looking strings keep their original formatting (matching
# every other quarter sheet, where these columns are stored as text).
$figures = $q1.Range("D2:G3")
$figures.NumberFormat = "@"

$q1.Range("D2").Value = "2.46"
$q1.Range("E2").Value = "81.85"
$q1.Range("F2").Value = "6.01"
$q1.Range("G2").Value = "0.1478"

$q1.Range("D3").Value = "2.46"
$q1.Range("E3").Value = "81.85"
$q1.Range("F3").Value = "6.01"
$q1.Range("G3").Value = "0.1478"

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(6)

# Shift existing rows 2-5 down to 3-6 (bottom-up so we never clobber a
# row before reading it).
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Range("A$dest").Value = $total.Range("A$r").Value2
    $total.Range("B$dest").Value = $total.Range("B$r").Value2
    $total.Range("C$dest").Value = $total.Range("C$r").Value2
    $total.Range("D$dest").Value = $total.Range("D$r").Value2
}

# Row 6 is brand new (sheet used to end at row 5) - give its index cell
# the same formatting as the rest of column A.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

# New first data row: 2022-Q1
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.3

# Renumber the running index column for every row.
for ($r = 3; $r -le 6; $r++) {
    $total.Range("A$r").Value = $r - 2
}

# Restore the original active sheet/selection.
$wb.Worksheets.Item(1).Select()
